$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.556.15"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.675.18"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'314.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.3949"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Value = "'0.3931"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "'1.393"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.45%  "
$ws.Range("D11").Value = "'50.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.04%  "
$ws.Range("D12").Value = "'0.08636"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").Value = "'7.286"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "'0.00001313"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "'7.646"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.51%  "
$ws.Range("D17").Value = "1.675.35"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "'93.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "'0.07015"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").Value = "'21.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'7.057"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("D24").Value = "24.552.17"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'2.347"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "'2.768"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.34%  "
$ws.Range("D27").Value = "'22.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").Value = "'5.832"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.30%  "
$ws.Range("D29").Value = "'158.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "'145.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "'8.261"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'2.542"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.15%  "
$ws.Range("D33").Value = "1.858.47"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("D34").Value = "'0.03068"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").Value = "'0.08254"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").Value = "'6.907"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.9935"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.53%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2795"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").Value = "'0.09613"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'1.512"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'10.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.86%  "
$ws.Range("D42").Value = "'0.7857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.64%  "
$ws.Range("E43").Value = "  -5.65%  "
$ws.Range("D44").Value = "'16.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.15%  "
$ws.Range("D45").Value = "'2.555"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("D46").Value = "'0.7077"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.74%  "
$ws.Range("D47").Value = "'4.169"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "'0.08639"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'1.325"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").Value = "'137.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
